$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update master data values as part of the 2nd May data refresh
$ws.Range("A3").Value = 10003
$ws.Range("A25").Value = 10003

# Scroll / selection changes: select full rows from 34 down to the end of the sheet
$ws.Rows("34:1048576").Select()
